$wb = $excel.ActiveWorkbook

# NOTE: worksheet lookup by name is case-insensitive, and this workbook has
# both "Vector_bf" (sheet 5) and "Vector_BF" (sheet 6) which would otherwise
# collide on the same match. Use 1-based sheet index instead, matching tab
# order: 1 Funciones_Objetivo, 2 Restricciones_del_lider,
# 3 Restricciones_del_follower, 4 Punto_modificado, 5 Vector_bf,
# 6 Vector_BF, 7 Vector_Alpha.

# ---------------------------------------------------------------------------
# Sheet: Restricciones_del_follower
#   Rows 2-6 hold (Expression, Function_Evaluation, Restriction_Set_Type,
#   Lambda_value, Beta_value, Gamma_value) records that were regenerated by
#   the experiment generator (alpha = 0 run). Format the range as Text first
#   so the numeric-looking literals round-trip as the exact original string
#   (e.g. "9.200000000000001") instead of being re-parsed/re-printed as a
#   binary double with different shortest-round-trip formatting.
# ---------------------------------------------------------------------------
$wsFollower = $wb.Worksheets.Item(3)
$wsFollower.Range("A2:F6").NumberFormat = "@"

$wsFollower.Range("A2").Value = "8.600000000000001 - y_1"
$wsFollower.Range("B2").Value = "-8.600000000000001"
$wsFollower.Range("C2").Value = "J_0_L0_v"
$wsFollower.Range("D2").Value = "0.75"
$wsFollower.Range("E2").Value = "6.7"
$wsFollower.Range("F2").Value = "9.200000000000001"

$wsFollower.Range("A3").Value = "-8.600000000000001 + y_1"
$wsFollower.Range("B3").Value = "4.600000000000001"
$wsFollower.Range("C3").Value = "J_0_L0_v"
$wsFollower.Range("D3").Value = "0.19"
$wsFollower.Range("E3").Value = "5.8"
$wsFollower.Range("F3").Value = "0.2"

$wsFollower.Range("A4").Value = "-5.000000000000002 - 2x + y_1 + 4y_2"
$wsFollower.Range("B4").Value = "-10.999999999999998"
$wsFollower.Range("C4").Value = "J_0_LP_v"
$wsFollower.Range("D4").Value = "0.24"
$wsFollower.Range("E4").Value = "5.8"
$wsFollower.Range("F4").Value = "8.0"

$wsFollower.Range("A5").Value = "-65.78 + 8x + y_1"
$wsFollower.Range("B5").Value = "17.4"
$wsFollower.Range("C5").Value = "J_Ne_L0_v"
$wsFollower.Range("D5").Value = "0.92"
$wsFollower.Range("E5").Value = "6.800000000000001"
$wsFollower.Range("F5").Value = "0.4"

$wsFollower.Range("A6").Value = "-7.400000000000002 - 2x - 2y_1"
$wsFollower.Range("B6").Value = "-19.400000000000002"
$wsFollower.Range("C6").Value = "J_Ne_L0_v"
$wsFollower.Range("D6").Value = "0.82"
$wsFollower.Range("E6").Value = "8.2"
$wsFollower.Range("F6").Value = "2.5"

# ---------------------------------------------------------------------------
# Sheet: Punto_modificado  (x, y_1, y_2 values)
# ---------------------------------------------------------------------------
$wsPunto = $wb.Worksheets.Item(4)
$wsPunto.Range("A2:C2").NumberFormat = "@"
$wsPunto.Range("A2").Value = "7.1"
$wsPunto.Range("B2").Value = "8.600000000000001"
$wsPunto.Range("C2").Value = "2.65"

# ---------------------------------------------------------------------------
# Sheet: Vector_bf
# ---------------------------------------------------------------------------
$wsBf = $wb.Worksheets.Item(5)
$wsBf.Range("A2:A3").NumberFormat = "@"
$wsBf.Range("A2").Value = "2.04"
$wsBf.Range("A3").Value = "-0.96"

# ---------------------------------------------------------------------------
# Sheet: Vector_BF
# ---------------------------------------------------------------------------
$wsBF = $wb.Worksheets.Item(6)
$wsBF.Range("A2:A3").NumberFormat = "@"
$wsBF.Range("A2").Value = "-25.400000000000006"
$wsBF.Range("A3").Value = "7.699999999999999"
# A4 ("-25.2") is unchanged.
